$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 55 with the same value as A54 (" 06-11-20")
$ws.Range("A55").Value = " 06-11-20"
